$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 5 for "Landesdatenbank NRW" (between "Open Data NRW" and "Open Data Moers") ---
$ws.Rows(5).Insert()

$ws.Range("A5").Value = "Landesdatenbank NRW"
$ws.Range("B5").Value = "Lokal"
$ws.Range("C5").Value = "https://www.landesdatenbank.nrw.de/"

# Give the new URL cell the same hyperlink look (blue underline) the other URL
# cells in column C already use.
$ws.Range("C5").Font.Underline = 1
$ws.Range("C5").Font.Color = 13391121

# The row-insert above does not shift the existing Hyperlinks collection, so the
# last data row (now row 13, "OUR World in Data") lost its working hyperlink --
# re-add it, which also restores a clickable link on the new text.
$ws.Hyperlinks.Add($ws.Range("C13"), "https://ourworldindata.org/") | Out-Null

# --- Widen column A to fit the longer labels ---
$ws.Columns("A").ColumnWidth = 39.3

# --- Page setup: fit to one page wide, centered, gridlines on print ---
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = $false
$ws.PageSetup.PrintGridlines = $true
$ws.PageSetup.CenterHorizontally = $true
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.PageSetup.HeaderMargin = 0.0
$ws.PageSetup.FooterMargin = 0.0
$ws.PageSetup.TopMargin = 54.0
$ws.PageSetup.BottomMargin = 54.0
$ws.PageSetup.LeftMargin = 50.4
$ws.PageSetup.RightMargin = 50.4

Write-Output "done"
